# LoginCredentialDetails.xlsx - "Added LogInPage Test Cases"
#
# Sheet1 gains 4 new LogIn test-case rows (APMS-T132..APMS-T136) and two of
# the existing rows get their TestCaseId changed from a numeric placeholder
# (290) to a proper text id. Row 5's Username/Password are also replaced.
# Finally, the active sheet/selection flips from sheet2 back to Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("verifyLogInCredentialsTest")

# --- Row 3 & 4: TestCaseId was a numeric placeholder (290); now real ids ---
$ws1.Range("A3").Value = "APMS-T132"
$ws1.Range("A4").Value = "APMS-T133"

# --- Row 5: new TestCaseId; Username/Password replaced later (see below) ---
$ws1.Range("A5").Value = "APMS-T134"

# --- Row 6 (new): APMS-T135 / Tesla / <blank, hyperlink-styled cell> ---
$ws1.Range("A6:B6").HorizontalAlignment = -4108   # xlCenter, matches rows above
$ws1.Range("A6").Value = "APMS-T135"
$ws1.Range("B6").Value = "Tesla"
$ws1.Range("C6").Style = "Hyperlink"

# --- Row 7 (new): APMS-T136 / <blank> / Tesla@123 (mailto hyperlink) ---
$ws1.Range("A7:B7").HorizontalAlignment = -4108   # xlCenter, matches rows above
$ws1.Range("A7").Value = "APMS-T136"
$ws1.Range("C7").Value = "Tesla@123"
[void]$ws1.Hyperlinks.Add($ws1.Range("C7"), "mailto:Tesla@123")

# --- Row 5 Username/Password replaced with new test values ---
# (done after rows 6/7 so the shared-string table order matches the recorded edit)
$ws1.Range("B5").Value = "Testing"
$ws1.Range("C5").Value = "Testing@123"

# --- View/selection: focus moves from sheet2 back to Sheet1 ---
[void]$ws2.Range("B2").Select()
$ws1.Activate()
[void]$ws1.Range("C8").Select()
